$d = $word.ActiveDocument

# 1) Insert two new, truly-empty paragraphs (no run) at the very start of
#    the document. Using InsertXML with a minimal package lets us put
#    bare `<w:p/>` elements in, matching the target markup exactly
#    (InsertParagraphBefore would instead leave a stray empty `<w:r/>`).
$startRange = $d.Range(0, 0)
$emptyParasXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$startRange.InsertXML($emptyParasXml) | Out-Null

# After the insert, paragraph numbering shifts by 2:
#  P1 = new empty
#  P2 = new empty
#  P3 = old "TTTclass:summary" + bookmark + "TTT"
#  P4 = old empty paragraph
#  P5 = old "TTTclass:scoreExampleTTT"
#  P6..P10 = remaining old paragraphs

# 2) Replace the text of paragraph 3 with the new wording (drops the
#    embedded _GoBack bookmark along with the old runs).
$p3 = $d.Paragraphs(3).Range
$p3.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark (wdCharacter = 1)
$p3.Text = "TTTclass:satisfactionTTT"

# 3) Clear paragraph 5's text and give it a `_GoBack` bookmark instead.
#    Adding a new `_GoBack` bookmark automatically removes the old one
#    Word kept around (bookmark names must be unique).
$p5 = $d.Paragraphs(5).Range
$p5.MoveEnd(1, -1) | Out-Null
$p5.Text = ""
$bmRange = $d.Paragraphs(5).Range
$bmRange.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 4) Remove the trailing paragraphs that are no longer needed (old
#    paragraphs 4 and 6-8, i.e. current paragraphs 6 through 10).
$deleteRange = $d.Range($d.Paragraphs(6).Range.Start, $d.Paragraphs(10).Range.End)
$deleteRange.Delete() | Out-Null
